$wb = $excel.ActiveWorkbook


# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 385.95
$ws.Range("I33").Value = 211.64516
$ws.Range("K33").Value = 211.64516
$ws.Range("M33").Value = 17.35484
$ws.Range("H64").Value = 5995.636
$ws.Range("J64").Value = 5549.8887
$ws.Range("L64").Value = 5549.8887
$ws.Range("N64").Value = -6045.8887
$ws.Range("H67").Value = 5995.636
$ws.Range("J67").Value = 5549.8887
$ws.Range("L67").Value = 5549.8887
$ws.Range("N67").Value = -7265.8887
$ws.Range("H106").Value = 3925041.2
$ws.Range("I106").Value = 4765221.5
$ws.Range("K106").Value = 4765221.5
$ws.Range("M106").Value = -4764590.5
$ws.Range("H112").Value = 2824.1025
$ws.Range("I112").Value = 898
$ws.Range("J112").Value = 3488.276
$ws.Range("K112").Value = 2694
$ws.Range("L112").Value = 10464.828
$ws.Range("M112").Value = -1586
$ws.Range("N112").Value = -12680.828

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2991.3408
$ws.Range("I61").Value = 2813.4866
$ws.Range("K61").Value = 2813.4866
$ws.Range("M61").Value = -2601.4866
$ws.Range("H63").Value = 159511
$ws.Range("I63").Value = 2500
$ws.Range("K63").Value = 2500
$ws.Range("M63").Value = -1814
$ws.Range("H66").Value = 159511
$ws.Range("I66").Value = 2500
$ws.Range("K66").Value = 12500
$ws.Range("M66").Value = -9068
$ws.Range("H74").Value = 1869.9
$ws.Range("I74").Value = 1809.7084
$ws.Range("J74").Value = 2110.6667
$ws.Range("K74").Value = 1809.7084
$ws.Range("L74").Value = 2110.6667
$ws.Range("M74").Value = -935.7084
$ws.Range("N74").Value = -3858.6667
$ws.Range("H77").Value = 1869.9
$ws.Range("I77").Value = 1809.7084
$ws.Range("J77").Value = 2110.6667
$ws.Range("K77").Value = 9048.542
$ws.Range("L77").Value = 10553.3335
$ws.Range("M77").Value = -4680.541999999999
$ws.Range("N77").Value = -19289.3335
$ws.Range("H110").Value = 2229.077
$ws.Range("I110").Value = 2226.0435
$ws.Range("K110").Value = 2226.0435
$ws.Range("M110").Value = -181.0435000000002
$ws.Range("N110").Value = -6342.3333
$ws.Range("H122").Value = 4661.222
$ws.Range("I122").Value = 3662.8518
$ws.Range("K122").Value = 10988.5554
$ws.Range("M122").Value = -8538.5554
$ws.Range("H132").Value = 3612.1777
$ws.Range("I132").Value = 3158.4583
$ws.Range("J132").Value = 4130.7144
$ws.Range("K132").Value = 9475.374899999999
$ws.Range("L132").Value = 12392.1432
$ws.Range("M132").Value = -6945.374899999999
$ws.Range("N132").Value = -17452.1432
$ws.Range("H136").Value = 2991.3408
$ws.Range("I136").Value = 2813.4866
$ws.Range("K136").Value = 8440.4598
$ws.Range("M136").Value = -5890.459800000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 345902.66
$ws.Range("I81").Value = 18854.5
$ws.Range("J81").Value = 999999
$ws.Range("K81").Value = 18854.5
$ws.Range("L81").Value = 999999
$ws.Range("M81").Value = -17793.5
$ws.Range("N81").Value = -1002121
$ws.Range("H84").Value = 345902.66
$ws.Range("I84").Value = 18854.5
$ws.Range("J84").Value = 999999
$ws.Range("K84").Value = 56563.5
$ws.Range("L84").Value = 2999997
$ws.Range("M84").Value = -51259.5
$ws.Range("N84").Value = -3010605
$ws.Range("H86").Value = 3709.7585
$ws.Range("I86").Value = 2554.45
$ws.Range("J86").Value = 6277.1113
$ws.Range("K86").Value = 2554.45
$ws.Range("L86").Value = 6277.1113
$ws.Range("M86").Value = -1431.45
$ws.Range("N86").Value = -8523.1113
$ws.Range("H89").Value = 3709.7585
$ws.Range("I89").Value = 2554.45
$ws.Range("J89").Value = 6277.1113
$ws.Range("K89").Value = 12772.25
$ws.Range("L89").Value = 31385.5565
$ws.Range("M89").Value = -7156.25
$ws.Range("N89").Value = -42617.5565
$ws.Range("H94").Value = 1114.4615
$ws.Range("I94").Value = 301.3
$ws.Range("J94").Value = 3825
$ws.Range("K94").Value = 301.3
$ws.Range("L94").Value = 3825
$ws.Range("M94").Value = 149.7
$ws.Range("N94").Value = -4727
$ws.Range("H105").Value = 2246.5
$ws.Range("I105").Value = 1828.6666
$ws.Range("K105").Value = 1828.6666
$ws.Range("M105").Value = -81.66660000000002
$ws.Range("H107").Value = 1925.6666
$ws.Range("I107").Value = 1202.3529
$ws.Range("K107").Value = 1202.3529
$ws.Range("M107").Value = 717.6470999999999
$ws.Range("H111").Value = 78990
$ws.Range("J111").Value = 78990
$ws.Range("L111").Value = 78990
$ws.Range("N111").Value = -87170
$ws.Range("H134").Value = 18523866
$ws.Range("I134").Value = 2774.875
$ws.Range("K134").Value = 8324.625
$ws.Range("M134").Value = -5789.625

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 147.875
$ws.Range("I7").Value = 147.875
$ws.Range("K7").Value = 147.875
$ws.Range("M7").Value = -34.875
$ws.Range("H88").Value = 41557.668
$ws.Range("J88").Value = 41557.668
$ws.Range("L88").Value = 41557.668
$ws.Range("N88").Value = -42369.668
$ws.Range("H91").Value = 41557.668
$ws.Range("J91").Value = 41557.668
$ws.Range("L91").Value = 41557.668
$ws.Range("N91").Value = -44365.668
$ws.Range("H92").Value = 64998.168
$ws.Range("J92").Value = 64998.168
$ws.Range("L92").Value = 64998.168
$ws.Range("N92").Value = -69990.168
$ws.Range("H105").Value = 2281.9167
$ws.Range("I105").Value = 2115
$ws.Range("J105").Value = 2782.6667
$ws.Range("K105").Value = 2115
$ws.Range("L105").Value = 2782.6667
$ws.Range("M105").Value = -368
$ws.Range("N105").Value = -6276.6667
$ws.Range("H137").Value = 118333.664
$ws.Range("J137").Value = 118333.664
$ws.Range("L137").Value = 118333.664
$ws.Range("N137").Value = -128533.664

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 821.5
$ws.Range("I55").Value = 821.5
$ws.Range("K55").Value = 2464.5
$ws.Range("M55").Value = -2287.5
$ws.Range("H68").Value = 909.6
$ws.Range("J68").Value = 999.6667
$ws.Range("L68").Value = 2999.0001
$ws.Range("N68").Value = -4621.0001
$ws.Range("H71").Value = 909.6
$ws.Range("J71").Value = 999.6667
$ws.Range("L71").Value = 8997.0003
$ws.Range("N71").Value = -17109.0003

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").Value = $null
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").Value = $null
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").Value = $null
$ws.Range("H97").Value = 1598.7
$ws.Range("I97").Value = 1061
$ws.Range("K97").Value = 1061
$ws.Range("M97").Value = -565
$ws.Range("H110").Value = 74008
$ws.Range("J110").Value = 74008
$ws.Range("L110").Value = 74008
$ws.Range("N110").Value = -82188
$ws.Range("H122").Value = 3243.4
$ws.Range("I122").Value = 2786.5454
$ws.Range("J122").Value = 4499.75
$ws.Range("K122").Value = 8359.6362
$ws.Range("L122").Value = 13499.25
$ws.Range("M122").Value = -5909.636200000001
$ws.Range("N122").Value = -18399.25
$ws.Range("H132").Value = 2500
$ws.Range("I132").Value = 2500
$ws.Range("K132").Value = 7500
$ws.Range("M132").Value = -4970

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H59").Value = 85475.29
$ws.Range("J59").Value = 85475.29
$ws.Range("L59").Value = 85475.29
$ws.Range("N59").Value = -86783.29
$ws.Range("H61").Value = 1751.4546
$ws.Range("I61").Value = 1612.25
$ws.Range("K61").Value = 1612.25
$ws.Range("M61").Value = -1410.25
$ws.Range("H93").Value = 967.45
$ws.Range("I93").Value = 540.2857
$ws.Range("K93").Value = 540.2857
$ws.Range("M93").Value = 707.7143
$ws.Range("H113").Value = 1751.4546
$ws.Range("I113").Value = 1612.25
$ws.Range("K113").Value = 1612.25
$ws.Range("M113").Value = 557.75
$ws.Range("H132").Value = 3613.8572
$ws.Range("I132").Value = 3519
$ws.Range("J132").Value = 3666.5557
$ws.Range("K132").Value = 10557
$ws.Range("L132").Value = 10999.6671
$ws.Range("M132").Value = -8027
$ws.Range("N132").Value = -16059.6671

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 13983.889
$ws.Range("J45").Value = 15481.875
$ws.Range("L45").Value = 15481.875
$ws.Range("N45").Value = -16463.875
$ws.Range("H100").Value = 2137.2
$ws.Range("I100").Value = 2116.6667
$ws.Range("K100").Value = 4233.3334
$ws.Range("M100").Value = -3692.3334
$ws.Range("H136").Value = 53464.25
$ws.Range("I136").Value = 2772.818
$ws.Range("K136").Value = 8318.454000000002
$ws.Range("M136").Value = -5768.454000000002
